# Scheduled runner update: refresh market-price columns (H-N) for a batch
# of leve rows across the crafting-class sheets (currentAveragePrice*,
# LevePrice*, LeveProfit* columns), as produced by the latest data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 651.5
$ws.Range("I28").Value = 324.46155
$ws.Range("J28").Value = 2068.6667
$ws.Range("K28").Value = 324.46155
$ws.Range("L28").Value = 2068.6667
$ws.Range("M28").Value = 160.53845
$ws.Range("N28").Value = -3038.6667

$ws.Range("H76").Value = 3209
$ws.Range("I76").Value = 3119.8
$ws.Range("J76").Value = 3283.3333
$ws.Range("K76").Value = 3119.8
$ws.Range("L76").Value = 3283.3333
$ws.Range("M76").Value = -2804.8
$ws.Range("N76").Value = -3913.3333

$ws.Range("H79").Value = 3209
$ws.Range("I79").Value = 3119.8
$ws.Range("J79").Value = 3283.3333
$ws.Range("K79").Value = 3119.8
$ws.Range("L79").Value = 3283.3333
$ws.Range("M79").Value = -2027.8
$ws.Range("N79").Value = -5467.3333

$ws.Range("H86").Value = 9010.5
$ws.Range("I86").Value = 1180
$ws.Range("J86").Value = 12925.75
$ws.Range("K86").Value = 1180
$ws.Range("L86").Value = 12925.75
$ws.Range("M86").Value = -57
$ws.Range("N86").Value = -15171.75

$ws.Range("H89").Value = 9010.5
$ws.Range("I89").Value = 1180
$ws.Range("J89").Value = 12925.75
$ws.Range("K89").Value = 5900
$ws.Range("L89").Value = 64628.75
$ws.Range("M89").Value = -284
$ws.Range("N89").Value = -75860.75

$ws.Range("H98").Value = 630.53845
$ws.Range("I98").Value = 309.7
$ws.Range("J98").Value = 1700
$ws.Range("K98").Value = 309.7
$ws.Range("L98").Value = 1700
$ws.Range("M98").Value = 1188.3
$ws.Range("N98").Value = -4696

$ws.Range("H106").Value = 8773693
$ws.Range("I106").Value = 11905868
$ws.Range("J106").Value = 3602.2
$ws.Range("K106").Value = 11905868
$ws.Range("L106").Value = 3602.2
$ws.Range("M106").Value = -11905237
$ws.Range("N106").Value = -4864.2

$ws.Range("H122").Value = 630.53845
$ws.Range("I122").Value = 309.7
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 929.0999999999999
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = 1520.9
$ws.Range("N122").Value = -10000

$ws.Range("H129").Value = 287950.16
$ws.Range("I129").Value = 495
$ws.Range("J129").Value = 305371.7
$ws.Range("K129").Value = 1485
$ws.Range("L129").Value = 916115.1000000001
$ws.Range("M129").Value = 3515
$ws.Range("N129").Value = -926115.1000000001

$ws.Range("H138").Value = 2647.2285
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2647.2285
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7941.685500000001
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -18221.6855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1529.1072
$ws.Range("I86").Value = 1400.75
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1400.75
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -277.75
$ws.Range("N86").Value = -4096

$ws.Range("H89").Value = 1529.1072
$ws.Range("I89").Value = 1400.75
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 7003.75
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -1387.75
$ws.Range("N89").Value = -20482

$ws.Range("H134").Value = 25513.174
$ws.Range("I134").Value = 28880.15
$ws.Range("J134").Value = 3066.6667
$ws.Range("K134").Value = 86640.45000000001
$ws.Range("L134").Value = 9200.000100000001
$ws.Range("M134").Value = -84105.45000000001
$ws.Range("N134").Value = -14270.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10606.081
$ws.Range("I31").Value = 11714.167
$ws.Range("J31").Value = 5857.143
$ws.Range("K31").Value = 11714.167
$ws.Range("L31").Value = 5857.143
$ws.Range("M31").Value = -11419.167
$ws.Range("N31").Value = -6447.143

$ws.Range("H34").Value = 10606.081
$ws.Range("I34").Value = 11714.167
$ws.Range("J34").Value = 5857.143
$ws.Range("K34").Value = 11714.167
$ws.Range("L34").Value = 5857.143
$ws.Range("M34").Value = -11512.167
$ws.Range("N34").Value = -6261.143

$ws.Range("H99").Value = 20837100
$ws.Range("I99").Value = 3336.842
$ws.Range("J99").Value = 100005400
$ws.Range("K99").Value = 3336.842
$ws.Range("L99").Value = 100005400
$ws.Range("M99").Value = -1838.842
$ws.Range("N99").Value = -100008396

$ws.Range("H126").Value = 20837100
$ws.Range("I126").Value = 3336.842
$ws.Range("J126").Value = 100005400
$ws.Range("K126").Value = 10010.526
$ws.Range("L126").Value = 300016200
$ws.Range("M126").Value = -7540.526
$ws.Range("N126").Value = -300021140

$ws.Range("H132").Value = 13313.955
$ws.Range("I132").Value = 15450.083
$ws.Range("J132").Value = 4769.4443
$ws.Range("K132").Value = 46350.249
$ws.Range("L132").Value = 14308.3329
$ws.Range("M132").Value = -43820.249
$ws.Range("N132").Value = -19368.3329

$ws.Range("H134").Value = 1024.5758
$ws.Range("I134").Value = 864.4091
$ws.Range("J134").Value = 1344.909
$ws.Range("K134").Value = 2593.2273
$ws.Range("L134").Value = 4034.727
$ws.Range("M134").Value = -58.22730000000001
$ws.Range("N134").Value = -9104.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5177
$ws.Range("I107").Value = 12640
$ws.Range("J107").Value = 201.66667
$ws.Range("K107").Value = 37920
$ws.Range("L107").Value = 605.00001
$ws.Range("M107").Value = -36000
$ws.Range("N107").Value = -4445.00001

$ws.Range("H131").Value = 782.05
$ws.Range("I131").Value = 286
$ws.Range("J131").Value = 808.1579
$ws.Range("K131").Value = 858
$ws.Range("L131").Value = 2424.4737
$ws.Range("M131").Value = 4182
$ws.Range("N131").Value = -12504.4737

$ws.Range("H132").Value = 1429.1428
$ws.Range("I132").Value = 999.5
$ws.Range("J132").Value = 1601
$ws.Range("K132").Value = 8995.5
$ws.Range("L132").Value = 14409
$ws.Range("M132").Value = -6465.5
$ws.Range("N132").Value = -19469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 55920
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 55920
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 55920
$ws.Range("N74").Value = -57792

$ws.Range("H77").Value = 55920
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 55920
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 167760
$ws.Range("N77").Value = -177120

$ws.Range("H80").Value = 3773.6843
$ws.Range("I80").Value = 3066.6667
$ws.Range("J80").Value = 4100
$ws.Range("K80").Value = 3066.6667
$ws.Range("L80").Value = 4100
$ws.Range("M80").Value = -2068.6667
$ws.Range("N80").Value = -6096

$ws.Range("H83").Value = 3773.6843
$ws.Range("I83").Value = 3066.6667
$ws.Range("J83").Value = 4100
$ws.Range("K83").Value = 15333.3335
$ws.Range("L83").Value = 20500
$ws.Range("M83").Value = -10341.3335
$ws.Range("N83").Value = -30484

$ws.Range("H132").Value = 38112.69
$ws.Range("I132").Value = 28664.078
$ws.Range("J132").Value = 127874.5
$ws.Range("K132").Value = 85992.234
$ws.Range("L132").Value = 383623.5
$ws.Range("M132").Value = -83462.234
$ws.Range("N132").Value = -388683.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1340.5869
$ws.Range("I132").Value = 1001.85
$ws.Range("J132").Value = 3598.8333
$ws.Range("K132").Value = 3005.55
$ws.Range("L132").Value = 10796.4999
$ws.Range("M132").Value = -475.5500000000002
$ws.Range("N132").Value = -15856.4999

$ws.Range("H136").Value = 15721.177
$ws.Range("I136").Value = 17503.5
$ws.Range("J136").Value = 2353.75
$ws.Range("K136").Value = 52510.5
$ws.Range("L136").Value = 7061.25
$ws.Range("M136").Value = -49960.5
$ws.Range("N136").Value = -12161.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1493.375
$ws.Range("I132").Value = 754.1905
$ws.Range("J132").Value = 6667.6665
$ws.Range("K132").Value = 2262.5715
$ws.Range("L132").Value = 20002.9995
$ws.Range("M132").Value = 267.4285
$ws.Range("N132").Value = -25062.9995

$ws.Range("H136").Value = 29413268
$ws.Range("I136").Value = 32259358
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 96778074
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -96775524
$ws.Range("N136").Value = -16100.0001
